$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 25.23990433333333
$ws.Range("H2").Value = 75.719713
$ws.Range("I2").Value = 0.05173702626903214
$ws.Range("J2").Value = 0.05173702626903214
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1575256666666667
$ws.Range("N2").Value = 0.472577
$ws.Range("O2").Value = 0.6985926944284299
$ws.Range("P2").Value = 0.69859269442843
$ws.Range("Q2").Value = 3.975932756711222
$ws.Range("R2").Value = 35.783394810401
$ws.Range("S2").Value = 0.03614310858299762
$ws.Range("T2").Value = 0.03614310858299762

$ws.Range("G3").Value = 25.23990433333333
$ws.Range("H3").Value = 75.719713
$ws.Range("I3").Value = 0.05173702626903214
$ws.Range("J3").Value = 0.05173702626903214
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.06796433333333333
$ws.Range("N3").Value = 0.203893
$ws.Range("O3").Value = 0.30140730557157
$ws.Range("P3").Value = 0.30140730557157
$ws.Range("Q3").Value = 1.715413271412111
$ws.Range("R3").Value = 15.438719442709
$ws.Range("S3").Value = 0.01559391768603451
$ws.Range("T3").Value = 0.01559391768603452

$ws.Range("I4").Value = 0.8454897015965644
$ws.Range("J4").Value = 0.8454897015965646
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1575256666666667
$ws.Range("N4").Value = 0.472577
$ws.Range("O4").Value = 0.6985926944284299
$ws.Range("P4").Value = 0.69859269442843
$ws.Range("Q4").Value = 64.97494043355778
$ws.Range("R4").Value = 584.77446390202
$ws.Range("S4").Value = 0.5906529287498331
$ws.Range("T4").Value = 0.5906529287498332

$ws.Range("I5").Value = 0.8454897015965644
$ws.Range("J5").Value = 0.8454897015965646
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.06796433333333333
$ws.Range("N5").Value = 0.203893
$ws.Range("O5").Value = 0.30140730557157
$ws.Range("P5").Value = 0.30140730557157
$ws.Range("Q5").Value = 28.03339038890889
$ws.Range("R5").Value = 252.30051350018
$ws.Range("S5").Value = 0.2548367728467313
$ws.Range("T5").Value = 0.2548367728467313

$ws.Range("G6").Value = 50.137933
$ws.Range("H6").Value = 150.413799
$ws.Range("I6").Value = 0.1027732721344034
$ws.Range("J6").Value = 0.1027732721344034
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1575256666666667
$ws.Range("N6").Value = 0.472577
$ws.Range("O6").Value = 0.6985926944284299
$ws.Range("P6").Value = 0.69859269442843
$ws.Range("Q6").Value = 7.898011321113667
$ws.Range("R6").Value = 71.082101890023
$ws.Range("S6").Value = 0.07179665709559911
$ws.Range("T6").Value = 0.07179665709559913

$ws.Range("G7").Value = 50.137933
$ws.Range("H7").Value = 150.413799
$ws.Range("I7").Value = 0.1027732721344034
$ws.Range("J7").Value = 0.1027732721344034
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.06796433333333333
$ws.Range("N7").Value = 0.203893
$ws.Range("O7").Value = 0.30140730557157
$ws.Range("P7").Value = 0.30140730557157
$ws.Range("Q7").Value = 3.407591191056333
$ws.Range("R7").Value = 30.66832071950699
$ws.Range("S7").Value = 0.03097661503880424
$ws.Range("T7").Value = 0.03097661503880424
